# Update DateEnd (column D) values for several ECV rows to reflect
# the newly included GIRAFE precipitation data coverage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D8"  = "2024-12-13"
    "D15" = "2024-02-29"
    "D18" = "2024-11-30"
    "D19" = "2023-12-31"
    "D20" = "2024-09-30"
    "D21" = "2024-12-02"
    "D22" = "2023-12-31"
}

foreach ($addr in $updates.Keys) {
    $rng = $ws.Range($addr)
    # Force text formatting first so Excel does not auto-convert the
    # date-like string into a serial date number, then write the value
    # and restore the cell's normal (unstyled) appearance so no visible
    # formatting change is introduced.
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$addr]
    $rng.Style = "Normal"
}
